# Fruta / hortaliza, semanal
# Insert 4 new weekly-report rows before the existing row 831, shifting the
# remaining historical rows (831-914) down to (835-918), then populate the
# 4 newly-inserted rows with the new period's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 831..834 (shifts old rows 831-914 down to 835-918).
$ws.Range("A831:A834").EntireRow.Insert()

# --- New row 831 ---
$ws.Cells.Item(831, 1).Value2 = 1
$ws.Cells.Item(831, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(831, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(831, 4).Value2 = 44946
$ws.Cells.Item(831, 5).Value2 = 15
$ws.Cells.Item(831, 6).Value2 = 100112024
$ws.Cells.Item(831, 7).Value2 = "Choclo"
$ws.Cells.Item(831, 8).Value2 = "Dulce o Americano"
$ws.Cells.Item(831, 9).Value2 = "Segunda"
$ws.Cells.Item(831, 10).Value2 = 160
$ws.Cells.Item(831, 11).Value2 = 11000
$ws.Cells.Item(831, 12).Value2 = 12000
$ws.Cells.Item(831, 13).Value2 = 11375
$ws.Cells.Item(831, 14).Value2 = "`$/malla 100 unidades"
$ws.Cells.Item(831, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(831, 16).Value2 = 1138
$ws.Cells.Item(831, 17).Value2 = 10
$ws.Cells.Item(831, 18).Value2 = "Hortaliza"

# --- New row 832 ---
$ws.Cells.Item(832, 1).Value2 = 1
$ws.Cells.Item(832, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(832, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(832, 4).Value2 = 44946
$ws.Cells.Item(832, 5).Value2 = 15
$ws.Cells.Item(832, 6).Value2 = 100112024
$ws.Cells.Item(832, 7).Value2 = "Choclo"
$ws.Cells.Item(832, 8).Value2 = "Lluteño"
$ws.Cells.Item(832, 9).Value2 = "Primera"
$ws.Cells.Item(832, 10).Value2 = 50
$ws.Cells.Item(832, 11).Value2 = 29000
$ws.Cells.Item(832, 12).Value2 = 30000
$ws.Cells.Item(832, 13).Value2 = 29600
$ws.Cells.Item(832, 14).Value2 = "`$/saco 50 unidades"
$ws.Cells.Item(832, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(832, 16).Value2 = 592
$ws.Cells.Item(832, 17).Value2 = 50
$ws.Cells.Item(832, 18).Value2 = "Hortaliza"

# --- New row 833 ---
$ws.Cells.Item(833, 1).Value2 = 1
$ws.Cells.Item(833, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(833, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(833, 4).Value2 = 44946
$ws.Cells.Item(833, 5).Value2 = 15
$ws.Cells.Item(833, 6).Value2 = 100112024
$ws.Cells.Item(833, 7).Value2 = "Choclo"
$ws.Cells.Item(833, 8).Value2 = "Lluteño"
$ws.Cells.Item(833, 9).Value2 = "Segunda"
$ws.Cells.Item(833, 10).Value2 = 25
$ws.Cells.Item(833, 11).Value2 = 24000
$ws.Cells.Item(833, 12).Value2 = 25000
$ws.Cells.Item(833, 13).Value2 = 24800
$ws.Cells.Item(833, 14).Value2 = "`$/saco 75 unidades"
$ws.Cells.Item(833, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(833, 16).Value2 = 331
$ws.Cells.Item(833, 17).Value2 = 75
$ws.Cells.Item(833, 18).Value2 = "Hortaliza"

# --- New row 834 ---
$ws.Cells.Item(834, 1).Value2 = 1
$ws.Cells.Item(834, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(834, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(834, 4).Value2 = 44946
$ws.Cells.Item(834, 5).Value2 = 15
$ws.Cells.Item(834, 6).Value2 = 100112024
$ws.Cells.Item(834, 7).Value2 = "Choclo"
$ws.Cells.Item(834, 8).Value2 = "Lluteño"
$ws.Cells.Item(834, 9).Value2 = "Tercera"
$ws.Cells.Item(834, 10).Value2 = 19
$ws.Cells.Item(834, 11).Value2 = 19000
$ws.Cells.Item(834, 12).Value2 = 20000
$ws.Cells.Item(834, 13).Value2 = 19579
$ws.Cells.Item(834, 14).Value2 = "`$/saco 100 unidades"
$ws.Cells.Item(834, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(834, 16).Value2 = 196
$ws.Cells.Item(834, 17).Value2 = 100
$ws.Cells.Item(834, 18).Value2 = "Hortaliza"
